# Reorders the "Recorded By" (column G) values so that the literal
# "System" token (exact case) - if present in the comma-separated list -
# is moved to the front of the list, ahead of any other recorder names.
#
# Example: "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
#          "system, backup@backdoor.com, System" -> "System, system, backup@backdoor.com"
#
# NOTE: this runtime's built-in string comparison operators (-eq, -ne,
# -ceq, -cne, -clike, -cmatch, ...) are all case-INSENSITIVE, so a
# manual, character-code based comparison is used to tell "System" apart
# from "system" / "SYSTEM" / etc.

function Test-ExactMatch($a, $b) {
    if ($a.Length -ne $b.Length) {
        return $false
    }
    for ($i = 0; $i -lt $a.Length; $i++) {
        $ca = [int][char]$a[$i]
        $cb = [int][char]$b[$i]
        if ($ca -ne $cb) {
            return $false
        }
    }
    return $true
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

# Column G is the 7th column ("Recorded By")
$col = 7

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $value = $cell.Value2

    if ($value -eq $null) {
        continue
    }

    $text = [string]$value
    if ($text -eq "") {
        continue
    }

    $parts = $text -split ","
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    $foundSystem = $false
    foreach ($p in $trimmed) {
        if (Test-ExactMatch $p "System") {
            $foundSystem = $true
        }
    }

    if ($foundSystem) {
        $rest = @()
        $systemRemoved = $false
        foreach ($p in $trimmed) {
            if ((-not $systemRemoved) -and (Test-ExactMatch $p "System")) {
                $systemRemoved = $true
            } else {
                $rest += $p
            }
        }
        $newParts = @("System") + $rest
        $newText = $newParts -join ", "

        if ($newText -ne $text) {
            $cell.Value = $newText
        }
    }
}
